# catching more jacks beta testing
# Adds new test-case rows for the "H991122AA" and "A041363AA" jack IDs to the
# Logic sheet's lookup/testing table and to the Sheet2 results table, and
# corrects the A041363AA-VD row in Sheet2 to reflect the new expected result.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Logic" sheet: append four new rows (23-26) describing the new jack
# types being tested, and update the selection/used range to match.
# ---------------------------------------------------------------------
$wsLogic = $wb.Worksheets.Item("Logic")

$wsLogic.Range("A23").Value = "H991122AA"
$wsLogic.Range("A24").Value = "H991122AA-d"
$wsLogic.Range("A25").Value = "A041363AA-D"
$wsLogic.Range("A26").Value = "A041363AA"

$wsLogic.Range("A2:A26").Select()

# ---------------------------------------------------------------------
# "Sheet2" sheet: fix row 10 (A041363AA-VD now resolves to CENET with a
# PG number instead of needing a manual update), then append the new
# generated/actual AP type rows (23-26) for the new jack types, plus six
# blank placeholder rows (27-32) matching the existing blank-row styling
# used elsewhere in column B.
# ---------------------------------------------------------------------
$wsSheet2 = $wb.Worksheets.Item("Sheet2")

$wsSheet2.Range("A10").Value = "A041363AA-VD, PG103139"
$wsSheet2.Range("B10").Value = "CENET"

$wsSheet2.Range("A23").Value = "H991122AA-DW, PG112220"
$wsSheet2.Range("B23").Value = "WIRELESS-AP-HOSPITALITY"

$wsSheet2.Range("A24").Value = "H991122AA-DW, PG112220"
$wsSheet2.Range("B24").Value = "WIRELESS-AP-HOSPITALITY"

$wsSheet2.Range("A25").Value = "A041363AA-D, PG103139"
$wsSheet2.Range("B25").Value = "CENET"

$wsSheet2.Range("A26").Value = "A041363AA"
$wsSheet2.Range("B26").Value = "Manual Update needed"

$wsSheet2.Range("B27").Value = ""
$wsSheet2.Range("B28").Value = ""
$wsSheet2.Range("B29").Value = ""
$wsSheet2.Range("B30").Value = ""
$wsSheet2.Range("B31").Value = ""
$wsSheet2.Range("B32").Value = ""
$wsSheet2.Range("B27:B32").Style = $wsSheet2.Range("B3").Style

$wsSheet2.Activate()
$wsSheet2.Range("B2:B26").Select()
